# Update "想去人数" (number of interested attendees) counts on the
# "展览" and "全部类型" sheets to reflect newly generated data
# (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F7").Value = 11319
    $ws.Range("F13").Value = 783
    $ws.Range("F14").Value = 12310
    $ws.Range("F15").Value = 12957
}
